$d = $word.ActiveDocument

# Replace each "{{" ... "}}" pair of 2-character runs with single-character
# "{" ... "}" runs, without letting the engine merge them into neighbouring
# runs. The engine coalesces adjacent runs that end up with identical
# formatting whenever a run's text is rewritten, so each target run is
# briefly given the opposite Bold value while its text is changed (which
# blocks the auto-merge because the neighbour's formatting no longer
# matches), and then restored to its original Bold value afterwards --
# right-hand run first, then the left-hand one, so the restore itself
# doesn't trigger a merge either.

$guard = 0
while ($true) {
    $guard = $guard + 1
    if ($guard -gt 50) { break }

    $rngOpen = $d.Content
    $findOpen = $rngOpen.Find
    $findOpen.ClearFormatting()
    $findOpen.Text = "{{"
    $findOpen.Forward = $true
    $findOpen.Wrap = 0
    $findOpen.MatchCase = $true
    $foundOpen = $findOpen.Execute()
    if (-not $foundOpen) { break }
    $s1 = $rngOpen.Start
    $e1 = $rngOpen.End

    $rngClose = $d.Content
    $findClose = $rngClose.Find
    $findClose.ClearFormatting()
    $findClose.Text = "}}"
    $findClose.Forward = $true
    $findClose.Wrap = 0
    $findClose.MatchCase = $true
    $foundClose = $findClose.Execute()
    if (-not $foundClose) { break }
    $s2 = $rngClose.Start
    $e2 = $rngClose.End

    $r1 = $d.Range($s1, $e1)
    $bold1 = $r1.Bold
    $r1.Bold = 1 - $bold1
    $r1.Text = "{"

    # The left-hand replacement shrank the document by (e1 - s1 - 1)
    # characters, so shift the right-hand pair's saved offsets accordingly.
    $shift = ($e1 - $s1) - 1
    $s2 = $s2 - $shift
    $e2 = $e2 - $shift

    $r2 = $d.Range($s2, $e2)
    $bold2 = $r2.Bold
    $r2.Bold = 1 - $bold2
    $r2.Text = "}"

    $r2restore = $d.Range($s2, $s2 + 1)
    $r2restore.Bold = $bold2

    $r1restore = $d.Range($s1, $s1 + 1)
    $r1restore.Bold = $bold1
}
